# Update "Đơn sale chính" sheet (sheet 1) and "Lương" sheet (sheet 3).
# Stops the empty placeholder report and fills in real per-service rows plus
# the salary calculation breakdown so downstream payroll code can consume it.

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: Don sale chinh ----
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(1,1).Value = 'Tiền tố'
$ws1.Cells.Item(1,2).Value = 'Mã dịch vụ'
$ws1.Cells.Item(1,3).Value = 'Ngày thực hiện'
$ws1.Cells.Item(1,4).Value = 'Cơ sở'
$ws1.Cells.Item(1,5).Value = 'Khách hàng'
$ws1.Cells.Item(1,6).Value = 'Nguồn khách'
$ws1.Cells.Item(1,7).Value = 'Nhóm dịch vụ'
$ws1.Cells.Item(1,8).Value = 'Tên dịch vụ'
$ws1.Cells.Item(1,9).Value = 'Sale chính'
$ws1.Cells.Item(1,10).Value = 'Đơn giá gốc'
$ws1.Cells.Item(1,11).Value = 'Sale phụ'
$ws1.Cells.Item(1,12).Value = 'Upsale'
$ws1.Cells.Item(1,13).Value = 'Đơn giá'
$ws1.Cells.Item(1,14).Value = 'Thanh toán lần đầu'
$ws1.Cells.Item(1,15).Value = 'Trả sau'
$ws1.Cells.Item(1,16).Value = 'Đã thanh toán'
$ws1.Cells.Item(1,17).Value = 'Dư nợ'
$ws1.Cells.Item(1,18).Value = 'Bác sĩ 1'
$ws1.Cells.Item(1,19).Value = 'Bác sĩ 2'
$ws1.Cells.Item(1,20).Value = 'Phụ phẫu 1'
$ws1.Cells.Item(1,21).Value = 'Phụ phẫu 2'
$ws1.Cells.Item(1,22).Value = 'Công phụ phẫu 1'
$ws1.Cells.Item(1,23).Value = 'Công phụ phẫu 2'
$ws1.Cells.Item(1,24).Value = 'Tỉ lệ chiết khấu sale chính'
$ws1.Cells.Item(1,25).Value = 'Tỉ lệ chiết khấu sale phụ'
$ws1.Cells.Item(1,26).Value = 'Chiết khấu sale chính'
$ws1.Cells.Item(1,27).Value = 'Chiết khấu sale phụ'
$ws1.Cells.Item(2,1).Value = 'HD-LUXURY'
$ws1.Cells.Item(2,2).Value = 507
$ws1.Cells.Item(2,3).Value = '07-01-2024'
$ws1.Cells.Item(2,4).Value = 'LONG XUYÊN'
$ws1.Cells.Item(2,5).Value = 'Nana'
$ws1.Cells.Item(2,6).Value = 'Khách cũ'
$ws1.Cells.Item(2,7).Value = 'Các ngoại khoa khác'
$ws1.Cells.Item(2,8).Value = 'Tiềm cằm'
$ws1.Cells.Item(2,9).Value = 'Nguyễn Phúc Nam'
$ws1.Cells.Item(2,10).Value = 3000000
$ws1.Cells.Item(2,13).Value = 3000000
$ws1.Cells.Item(2,14).Value = 3000000
$ws1.Cells.Item(2,15).Value = 0
$ws1.Cells.Item(2,16).Value = 3000000
$ws1.Cells.Item(2,17).Value = 0
$ws1.Cells.Item(2,18).Value = 'Đặng Ngọc Mai'
$ws1.Cells.Item(2,20).Value = 'Đào Vương Anh'
$ws1.Cells.Item(2,22).Value = 0
$ws1.Cells.Item(2,23).Value = 0
$ws1.Cells.Item(2,24).Value = 0.1
$ws1.Cells.Item(2,25).Value = 0
$ws1.Cells.Item(2,26).Value = 300000
$ws1.Cells.Item(2,27).Value = 0
$ws1.Cells.Item(3,1).Value = 'HD-LUXURY'
$ws1.Cells.Item(3,2).Value = 530
$ws1.Cells.Item(3,3).Value = '07-08-2024'
$ws1.Cells.Item(3,4).Value = 'LONG XUYÊN'
$ws1.Cells.Item(3,5).Value = 'Võ thị nga'
$ws1.Cells.Item(3,6).Value = 'Khách cũ'
$ws1.Cells.Item(3,7).Value = 'Vùng mắt'
$ws1.Cells.Item(3,8).Value = 'Cắt mí'
$ws1.Cells.Item(3,9).Value = 'Nguyễn Phúc Nam'
$ws1.Cells.Item(3,10).Value = 4000000
$ws1.Cells.Item(3,13).Value = 4000000
$ws1.Cells.Item(3,14).Value = 3000000
$ws1.Cells.Item(3,15).Value = 0
$ws1.Cells.Item(3,16).Value = 3000000
$ws1.Cells.Item(3,17).Value = 1000000
$ws1.Cells.Item(3,18).Value = 'Nguyễn Hoàng Yến Quyên'
$ws1.Cells.Item(3,20).Value = 'Đào Vương Anh'
$ws1.Cells.Item(3,22).Value = 50000
$ws1.Cells.Item(3,23).Value = 0
$ws1.Cells.Item(3,24).Value = 0.1
$ws1.Cells.Item(3,25).Value = 0
$ws1.Cells.Item(3,26).Value = 300000
$ws1.Cells.Item(3,27).Value = 0
$ws1.Cells.Item(4,1).Value = 'HD-LUXURY'
$ws1.Cells.Item(4,2).Value = 533
$ws1.Cells.Item(4,3).Value = '07-09-2024'
$ws1.Cells.Item(4,4).Value = 'LONG XUYÊN'
$ws1.Cells.Item(4,5).Value = 'Diễm hương'
$ws1.Cells.Item(4,6).Value = 'Cá nhân'
$ws1.Cells.Item(4,7).Value = 'Môi'
$ws1.Cells.Item(4,8).Value = 'Tiêm môi'
$ws1.Cells.Item(4,9).Value = 'Nguyễn Phúc Nam'
$ws1.Cells.Item(4,10).Value = 1350000
$ws1.Cells.Item(4,13).Value = 1350000
$ws1.Cells.Item(4,14).Value = 1350000
$ws1.Cells.Item(4,15).Value = 0
$ws1.Cells.Item(4,16).Value = 1350000
$ws1.Cells.Item(4,17).Value = 0
$ws1.Cells.Item(4,18).Value = 'Đặng Ngọc Mai'
$ws1.Cells.Item(4,20).Value = 'Đào Vương Anh'
$ws1.Cells.Item(4,22).Value = 0
$ws1.Cells.Item(4,23).Value = 0
$ws1.Cells.Item(4,24).Value = 0.1
$ws1.Cells.Item(4,25).Value = 0
$ws1.Cells.Item(4,26).Value = 135000
$ws1.Cells.Item(4,27).Value = 0
$ws1.Cells.Item(5,1).Value = 'Tổng'
$ws1.Cells.Item(5,2).Value = 3
$ws1.Cells.Item(5,10).Value = 8350000
$ws1.Cells.Item(5,12).Value = 0
$ws1.Cells.Item(5,13).Value = 8350000
$ws1.Cells.Item(5,14).Value = 7350000
$ws1.Cells.Item(5,15).Value = 0
$ws1.Cells.Item(5,16).Value = 7350000
$ws1.Cells.Item(5,17).Value = 1000000
$ws1.Cells.Item(5,22).Value = 50000
$ws1.Cells.Item(5,23).Value = 0
$ws1.Cells.Item(5,24).Value = 0.3
$ws1.Cells.Item(5,25).Value = 0
$ws1.Cells.Item(5,26).Value = 735000
$ws1.Cells.Item(5,27).Value = 0

# ---- Sheet 3: Luong ----
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(1,1).Value = 'Danh mục'
$ws3.Cells.Item(1,2).Value = 12
$ws3.Cells.Item(2,1).Value = 'Ngày công'
$ws3.Cells.Item(2,2).Value = 11.5
$ws3.Cells.Item(3,1).Value = 'Phụ cấp'
$ws3.Cells.Item(3,2).Value = 402500
$ws3.Cells.Item(4,1).Value = 'Lương cơ bản tại CẦN THƠ'
$ws3.Cells.Item(5,1).Value = 'Chiết khấu sale chính tại CẦN THƠ'
$ws3.Cells.Item(5,2).Value = 0
$ws3.Cells.Item(6,1).Value = 'Chiết khấu sale phụ tại CẦN THƠ'
$ws3.Cells.Item(6,2).Value = 0
$ws3.Cells.Item(7,1).Value = 'Đơn 1 bác sĩ tại CẦN THƠ'
$ws3.Cells.Item(7,2).Value = 0
$ws3.Cells.Item(8,1).Value = 'Đơn 2 bác sĩ tại CẦN THƠ'
$ws3.Cells.Item(8,2).Value = 0
$ws3.Cells.Item(9,1).Value = 'Công phụ phẫu 1 tại CẦN THƠ'
$ws3.Cells.Item(9,2).Value = 0
$ws3.Cells.Item(10,1).Value = 'Công phụ phẫu 2 tại CẦN THƠ'
$ws3.Cells.Item(10,2).Value = 0
$ws3.Cells.Item(11,1).Value = 'Lương cơ bản tại LONG XUYÊN'
$ws3.Cells.Item(11,2).Value = 3482857.142857143
$ws3.Cells.Item(12,1).Value = 'Chiết khấu sale chính tại LONG XUYÊN'
$ws3.Cells.Item(12,2).Value = 735000
$ws3.Cells.Item(13,1).Value = 'Chiết khấu sale phụ tại LONG XUYÊN'
$ws3.Cells.Item(13,2).Value = 0
$ws3.Cells.Item(14,1).Value = 'Đơn 1 bác sĩ tại LONG XUYÊN'
$ws3.Cells.Item(14,2).Value = 0
$ws3.Cells.Item(15,1).Value = 'Đơn 2 bác sĩ tại LONG XUYÊN'
$ws3.Cells.Item(15,2).Value = 0
$ws3.Cells.Item(16,1).Value = 'Công phụ phẫu 1 tại LONG XUYÊN'
$ws3.Cells.Item(16,2).Value = 0
$ws3.Cells.Item(17,1).Value = 'Công phụ phẫu 2 tại LONG XUYÊN'
$ws3.Cells.Item(17,2).Value = 0
$ws3.Cells.Item(18,1).Value = 'Lương cơ bản tại SÓC TRĂNG'
$ws3.Cells.Item(19,1).Value = 'Chiết khấu sale chính tại SÓC TRĂNG'
$ws3.Cells.Item(19,2).Value = 0
$ws3.Cells.Item(20,1).Value = 'Chiết khấu sale phụ tại SÓC TRĂNG'
$ws3.Cells.Item(20,2).Value = 0
$ws3.Cells.Item(21,1).Value = 'Đơn 1 bác sĩ tại SÓC TRĂNG'
$ws3.Cells.Item(21,2).Value = 0
$ws3.Cells.Item(22,1).Value = 'Đơn 2 bác sĩ tại SÓC TRĂNG'
$ws3.Cells.Item(22,2).Value = 0
$ws3.Cells.Item(23,1).Value = 'Công phụ phẫu 1 tại SÓC TRĂNG'
$ws3.Cells.Item(23,2).Value = 0
$ws3.Cells.Item(24,1).Value = 'Công phụ phẫu 2 tại SÓC TRĂNG'
$ws3.Cells.Item(24,2).Value = 0
